$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '31.218.00'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.54%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.998.48'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +6.63%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7836'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +65.84%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '257.07'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +4.54%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9994'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3501'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +21.95%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '28.86'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +31.64%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07035'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +8.12%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8588'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +17.74%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08211'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +5.16%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.998.94'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +6.67%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '100.96'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.29%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.605'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +8.41%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.71'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +19.99%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '273.68'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -3.80%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '31.228.59'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +2.66%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.961'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +11.84%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007925'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +5.81%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.263.49'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +7.01%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9991'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.04%  '

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.12%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.126'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +12.44%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.05'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +11.02%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.88'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.85%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1480'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +52.74%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.98'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +5.13%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.360'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +24.33%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.611'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +8.02%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.621'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +9.16%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.354'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.33%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.456'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +7.23%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05218'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +8.23%  '

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +9.16%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7759'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +12.15%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.789'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.65%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02003'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +5.39%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.922'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.95%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.740'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +6.89%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '79.67'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +4.82%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.153'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +9.86%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4690'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +11.07%  '

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +5.46%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8474'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +2.56%  '

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.01%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.743'

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.931'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.79%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4331'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +10.29%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '36.88'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +5.33%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.514'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +13.75%  '
